$d = $word.ActiveDocument

$replacements = @(
    @("46×20=920", "69×40=2760"),
    @("19×48=912", "77×32=2464"),
    @("34×28=952", "86×24=2064"),
    @("37×68=2516", "87×12=1044"),
    @("63×73=4599", "53×81=4293"),
    @("37×64=2368", "66×88=5808"),
    @("97×17=1649", "67×43=2881"),
    @("75×35=2625", "58×94=5452"),
    @("89×21=1869", "59×71=4189"),
    @("83×69=5727", "31×68=2108"),
    @("40×77=3080", "16×36=576"),
    @("84×14=1176", "18×44=792"),
    @("67×95=6365", "27×25=675"),
    @("31×74=2294", "82×31=2542"),
    @("72×56=4032", "46×85=3910"),
    @("11×57=627", "75×17=1275"),
    @("24×65=1560", "66×71=4686"),
    @("69×51=3519", "81×71=5751"),
    @("50×14=700", "20×52=1040"),
    @("64×71=4544", "27×49=1323"),
    @("68×94=6392", "56×66=3696"),
    @("87×89=7743", "83×87=7221"),
    @("25×90=2250", "33×16=528"),
    @("63×79=4977", "50×38=1900"),
    @("70×80=5600", "72×36=2592")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
